# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" (Exhibition), "演出" (Performance) and "全部类型" (All types) sheets,
# reflecting the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 4589
$wsExpo.Range("F9").Value = 3085
$wsExpo.Range("F13").Value = 606
$wsExpo.Range("F19").Value = 1316
$wsExpo.Range("F20").Value = 119
$wsExpo.Range("F31").Value = 3668
$wsExpo.Range("F34").Value = 391
$wsExpo.Range("F36").Value = 1756

# --- 演出 (Performance) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 22
$wsShow.Range("F3").Value = 39

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 4589
$wsAll.Range("F9").Value = 3085
$wsAll.Range("F13").Value = 606
$wsAll.Range("F16").Value = 22
$wsAll.Range("F20").Value = 1316
$wsAll.Range("F21").Value = 119
$wsAll.Range("F32").Value = 3668
$wsAll.Range("F33").Value = 39
$wsAll.Range("F36").Value = 391
$wsAll.Range("F38").Value = 1756
